# Auto-generated Excel COM-interop script to apply numeric updates
# described by the xml diff (Sheets/Bahamut_Profits.xlsx -> workbook sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1533.5454
$ws.Range("I15").Value = 1533.5454
$ws.Range("K15").Value = 4600.6362
$ws.Range("M15").Value = -4431.6362
$ws.Range("H116").Value = 2545
$ws.Range("I116").Value = 2810.8333
$ws.Range("J116").Value = 1861.4286
$ws.Range("K116").Value = 2810.8333
$ws.Range("L116").Value = 1861.4286
$ws.Range("M116").Value = 631.1667000000002
$ws.Range("N116").Value = -8745.428599999999
$ws.Range("H132").Value = 1724.5333
$ws.Range("I132").Value = 1360.1143
$ws.Range("K132").Value = 4080.3429
$ws.Range("M132").Value = -1550.3429
$ws.Range("H137").Value = 8929848
$ws.Range("I137").Value = 1240.9791
$ws.Range("J137").Value = 62501490
$ws.Range("K137").Value = 3722.9373
$ws.Range("L137").Value = 187504470
$ws.Range("M137").Value = -1172.9373
$ws.Range("N137").Value = -187509570
$ws.Range("H141").Value = 1251.4
$ws.Range("I141").Value = 797.36365
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 2392.09095
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = 2787.90905
$ws.Range("N141").Value = -17860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 14200
$ws.Range("I21").Value = 3000
$ws.Range("J21").Value = 19800
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 19800
$ws.Range("M21").Value = -2626
$ws.Range("N21").Value = -20548
$ws.Range("H61").Value = 1888.0952
$ws.Range("I61").Value = 1683.6
$ws.Range("J61").Value = 2074
$ws.Range("K61").Value = 1683.6
$ws.Range("L61").Value = 2074
$ws.Range("M61").Value = -1471.6
$ws.Range("N61").Value = -2498
$ws.Range("H74").Value = 851.44183
$ws.Range("I74").Value = 857.2
$ws.Range("K74").Value = 857.2
$ws.Range("M74").Value = 16.79999999999995
$ws.Range("H77").Value = 851.44183
$ws.Range("I77").Value = 857.2
$ws.Range("K77").Value = 4286
$ws.Range("M77").Value = 82
$ws.Range("H136").Value = 1888.0952
$ws.Range("I136").Value = 1683.6
$ws.Range("J136").Value = 2074
$ws.Range("K136").Value = 5050.799999999999
$ws.Range("L136").Value = 6222
$ws.Range("M136").Value = -2500.799999999999
$ws.Range("N136").Value = -11322

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 684
$ws.Range("I5").Value = 351.33334
$ws.Range("J5").Value = 1016.6667
$ws.Range("K5").Value = 351.33334
$ws.Range("L5").Value = 1016.6667
$ws.Range("M5").Value = -238.33334
$ws.Range("N5").Value = -1242.6667
$ws.Range("H99").Value = 83335690
$ws.Range("I99").Value = 100002100
$ws.Range("J99").Value = 3650
$ws.Range("K99").Value = 100002100
$ws.Range("L99").Value = 3650
$ws.Range("M99").Value = -100000602
$ws.Range("N99").Value = -6646
$ws.Range("H129").Value = 49989.5
$ws.Range("J129").Value = 49989.5
$ws.Range("L129").Value = 49989.5
$ws.Range("N129").Value = -59989.5
$ws.Range("H130").Value = 39903.332
$ws.Range("J130").Value = 39903.332
$ws.Range("L130").Value = 39903.332
$ws.Range("N130").Value = -49943.332
$ws.Range("H134").Value = 42261.92
$ws.Range("I134").Value = 2810.5715
$ws.Range("J134").Value = 92472.73
$ws.Range("K134").Value = 8431.7145
$ws.Range("L134").Value = 277418.19
$ws.Range("M134").Value = -5896.7145
$ws.Range("N134").Value = -282488.19

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1350.0968
$ws.Range("I58").Value = 990.2222
$ws.Range("J58").Value = 1497.3182
$ws.Range("K58").Value = 990.2222
$ws.Range("L58").Value = 1497.3182
$ws.Range("M58").Value = -787.2222
$ws.Range("N58").Value = -1903.3182
$ws.Range("H134").Value = 1955.5
$ws.Range("I134").Value = 1301.875
$ws.Range("J134").Value = 4570
$ws.Range("K134").Value = 3905.625
$ws.Range("L134").Value = 13710
$ws.Range("M134").Value = -1370.625
$ws.Range("N134").Value = -18780
$ws.Range("H136").Value = 1350.0968
$ws.Range("I136").Value = 990.2222
$ws.Range("J136").Value = 1497.3182
$ws.Range("K136").Value = 2970.6666
$ws.Range("L136").Value = 4491.9546
$ws.Range("M136").Value = -420.6666
$ws.Range("N136").Value = -9591.954600000001
$ws.Range("H137").Value = 37390
$ws.Range("I137").Value = 37000
$ws.Range("J137").Value = 37780
$ws.Range("K137").Value = 37000
$ws.Range("L137").Value = 37780
$ws.Range("M137").Value = -31900
$ws.Range("N137").Value = -47980

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1394.527
$ws.Range("I68").Value = 786.28
$ws.Range("J68").Value = 1704.8572
$ws.Range("K68").Value = 2358.84
$ws.Range("L68").Value = 5114.571599999999
$ws.Range("M68").Value = -1547.84
$ws.Range("N68").Value = -6736.571599999999
$ws.Range("H71").Value = 1394.527
$ws.Range("I71").Value = 786.28
$ws.Range("J71").Value = 1704.8572
$ws.Range("K71").Value = 7076.52
$ws.Range("L71").Value = 15343.7148
$ws.Range("M71").Value = -3020.52
$ws.Range("N71").Value = -23455.7148
$ws.Range("H129").Value = 43111.5
$ws.Range("I129").Value = 1143.2222
$ws.Range("J129").Value = 68292.47
$ws.Range("K129").Value = 3429.6666
$ws.Range("L129").Value = 204877.41
$ws.Range("M129").Value = 1570.3334
$ws.Range("N129").Value = -214877.41
$ws.Range("H131").Value = 44918
$ws.Range("I131").Value = 112758.89
$ws.Range("J131").Value = 4213.467
$ws.Range("K131").Value = 338276.67
$ws.Range("L131").Value = 12640.401
$ws.Range("M131").Value = -333236.67
$ws.Range("N131").Value = -22720.401

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 61.166668
$ws.Range("I2").Value = 33.083332
$ws.Range("J2").Value = 117.333336
$ws.Range("K2").Value = 33.083332
$ws.Range("L2").Value = 117.333336
$ws.Range("M2").Value = 79.916668
$ws.Range("N2").Value = -343.333336
$ws.Range("H123").Value = 28081.5
$ws.Range("J123").Value = 28081.5
$ws.Range("L123").Value = 28081.5
$ws.Range("N123").Value = -32981.5
$ws.Range("H130").Value = 43320
$ws.Range("J130").Value = 43320
$ws.Range("L130").Value = 43320
$ws.Range("N130").Value = -53360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1200.2858
$ws.Range("I16").Value = 1200.2858
$ws.Range("K16").Value = 1200.2858
$ws.Range("M16").Value = -1030.2858
$ws.Range("H61").Value = 734.53845
$ws.Range("I61").Value = 750.0909
$ws.Range("J61").Value = 649
$ws.Range("K61").Value = 750.0909
$ws.Range("L61").Value = 649
$ws.Range("M61").Value = -548.0909
$ws.Range("N61").Value = -1053
$ws.Range("H113").Value = 734.53845
$ws.Range("I113").Value = 750.0909
$ws.Range("J113").Value = 649
$ws.Range("K113").Value = 750.0909
$ws.Range("L113").Value = 649
$ws.Range("M113").Value = 1419.9091
$ws.Range("N113").Value = -4989
$ws.Range("H132").Value = 1669126.5
$ws.Range("I132").Value = 2254535.8
$ws.Range("J132").Value = 2961.3076
$ws.Range("K132").Value = 6763607.399999999
$ws.Range("L132").Value = 8883.9228
$ws.Range("M132").Value = -6761077.399999999
$ws.Range("N132").Value = -13943.9228
$ws.Range("H136").Value = 2016.0392
$ws.Range("I136").Value = 1205.8636
$ws.Range("K136").Value = 3617.5908
$ws.Range("M136").Value = -1067.5908
$ws.Range("H137").Value = 70000
$ws.Range("J137").Value = 70000
$ws.Range("L137").Value = 70000
$ws.Range("N137").Value = -80200
